# Updated temp sensor purchase link.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update unit price for the TMP36 temperature sensor line (row 18)
$ws.Range("E18").Value = 1.42

# Total cost becomes a computed formula instead of a static number
$ws.Range("F18").Formula = "=D18*E18"

# Point the "Source (link)" cell at the new DigiKey product page and
# refresh the displayed text to match the new URL
$ws.Range("G18").Value = "http://www.digikey.com/product-detail/en/TMP36GT9Z/TMP36GT9Z-ND/820404"
$ws.Hyperlinks.Add($ws.Range("G18"), "http://www.digikey.com/product-detail/en/TMP36GT9Z/TMP36GT9Z-ND/820404")

# Recalculate dependent totals
$excel.CalculateFull()

# Leave selection where Excel would land after editing G18 (one row below)
$ws.Range("F19").Select()
